# Updates the cryptocurrency price/volume snapshot (rows 2-51) to the
# refreshed values from the latest GitHub Actions data pull.
#
# Column D (Price) and column E (Volume 1h) hold text, not numbers
# (e.g. "25.303.17", "0.9980", "  -2.33%  "). Several of the Price values
# parse as plain numbers (e.g. "0.9980", "1.216"), so assigning them as a
# bare string would let Excel's COM layer auto-coerce them into numeric
# values and strip the meaningful trailing zeros. Prefixing the literal
# with a leading apostrophe (an extra leading single-quote inside the
# PowerShell single-quoted string, i.e. '''...' ) forces Excel to keep the
# text verbatim, matching how the source data is stored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.303.17'
$ws.Range("E2").Value = '''  -2.33%  '

$ws.Range("D3").Value = '''1.661.78'
$ws.Range("E3").Value = '''  -4.07%  '

$ws.Range("D4").Value = '''0.9966'
$ws.Range("E4").Value = '''  -0.30%  '

$ws.Range("D5").Value = '''235.77'
$ws.Range("E5").Value = '''  -4.17%  '

$ws.Range("D6").Value = '''0.9976'
$ws.Range("E6").Value = '''  -0.24%  '

$ws.Range("D7").Value = '''0.4805'
$ws.Range("E7").Value = '''  -4.62%  '

$ws.Range("D8").Value = '''0.2595'
$ws.Range("E8").Value = '''  -4.74%  '

$ws.Range("D9").Value = '''0.06145'
$ws.Range("E9").Value = '''  -0.42%  '

$ws.Range("D10").Value = '''0.07072'
$ws.Range("E10").Value = '''  -2.34%  '

$ws.Range("D11").Value = '''1.650.58'
$ws.Range("E11").Value = '''  -4.79%  '

$ws.Range("D12").Value = '''14.66'
$ws.Range("E12").Value = '''  -3.25%  '

$ws.Range("D13").Value = '''0.5854'
$ws.Range("E13").Value = '''  -10.52%  '

$ws.Range("D14").Value = '''4.362'
$ws.Range("E14").Value = '''  -8.69%  '

$ws.Range("D15").Value = '''74.31'
$ws.Range("E15").Value = '''  -3.62%  '

$ws.Range("D16").Value = '''0.9980'
$ws.Range("E16").Value = '''  -0.20%  '

$ws.Range("D17").Value = '''0.9975'
$ws.Range("E17").Value = '''  -0.18%  '

$ws.Range("D18").Value = '''25.275.53'
$ws.Range("E18").Value = '''  -2.49%  '

$ws.Range("D19").Value = '''0.000006691'
$ws.Range("E19").Value = '''  -1.78%  '

$ws.Range("D20").Value = '''11.39'
$ws.Range("E20").Value = '''  -4.01%  '

$ws.Range("D21").Value = '''1.859.82'
$ws.Range("E21").Value = '''  -5.15%  '

$ws.Range("D22").Value = '''4.365'
$ws.Range("E22").Value = '''  -4.93%  '

$ws.Range("D23").Value = '''8.599'
$ws.Range("E23").Value = '''  -2.19%  '

$ws.Range("D24").Value = '''5.315'
$ws.Range("E24").Value = '''  -3.06%  '

$ws.Range("D25").Value = '''134.54'
$ws.Range("E25").Value = '''  +0.44%  '

$ws.Range("D26").Value = '''15.11'
$ws.Range("E26").Value = '''  -0.85%  '

$ws.Range("D27").Value = '''1.379'
$ws.Range("E27").Value = '''  -3.35%  '

$ws.Range("D28").Value = '''104.61'
$ws.Range("E28").Value = '''  -0.45%  '

$ws.Range("D29").Value = '''1.677'
$ws.Range("E29").Value = '''  -6.31%  '

$ws.Range("D30").Value = '''3.953'
$ws.Range("E30").Value = '''  -1.12%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''3.611'
$ws.Range("E31").Value = '''  -2.30%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.07644'
$ws.Range("E32").Value = '''  -5.83%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.04337'
$ws.Range("E33").Value = '''  -8.22%  '

$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '''0.9975'
$ws.Range("E34").Value = '''  -0.15%  '

$ws.Range("E35").Value = '''  -2.21%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.6012'
$ws.Range("E36").Value = '''  -1.86%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '''0.9419'
$ws.Range("E37").Value = '''  -5.57%  '

$ws.Range("D38").Value = '''2.609'
$ws.Range("E38").Value = '''  -5.01%  '

$ws.Range("D39").Value = '''0.8453'
$ws.Range("E39").Value = '''  -4.81%  '

$ws.Range("D40").Value = '''0.9981'
$ws.Range("E40").Value = '''  -0.15%  '

$ws.Range("D41").Value = '''0.01497'
$ws.Range("E41").Value = '''  -6.56%  '

$ws.Range("D42").Value = '''99.32'
$ws.Range("E42").Value = '''  -1.38%  '

$ws.Range("D43").Value = '''1.809'
$ws.Range("E43").Value = '''  -7.29%  '

$ws.Range("D44").Value = '''0.3735'
$ws.Range("E44").Value = '''  -4.60%  '

$ws.Range("D45").Value = '''4.657'
$ws.Range("E45").Value = '''  -7.27%  '

$ws.Range("D46").Value = '''6.203'
$ws.Range("E46").Value = '''  -1.96%  '

$ws.Range("D47").Value = '''0.1111'
$ws.Range("E47").Value = '''  -5.67%  '

$ws.Range("D48").Value = '''0.05247'
$ws.Range("E48").Value = '''  -0.57%  '

$ws.Range("D49").Value = '''29.42'
$ws.Range("E49").Value = '''  -4.26%  '

$ws.Range("D50").Value = '''1.216'
$ws.Range("E50").Value = '''  -1.54%  '

$ws.Range("D51").Value = '''0.9990'
$ws.Range("E51").Value = '''  -0.23%  '
